# handled the time sync issues
# Update the Start time / End time / Time taken columns (E:G) for the
# 5 data rows with refreshed timestamps from the latest test run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2023-07-11 08:44:18"
$ws.Range("F2").Value = "2023-07-11 08:48:13"
$ws.Range("G2").Value = "00:03:55"

$ws.Range("E3").Value = "2023-07-11 08:48:15"
$ws.Range("F3").Value = "2023-07-11 08:52:10"
$ws.Range("G3").Value = "00:03:55"

$ws.Range("E4").Value = "2023-07-11 08:52:11"
$ws.Range("F4").Value = "2023-07-11 08:56:06"
$ws.Range("G4").Value = "00:03:55"

$ws.Range("E5").Value = "2023-07-11 08:56:08"
$ws.Range("F5").Value = "2023-07-11 09:00:02"
$ws.Range("G5").Value = "00:03:54"

$ws.Range("E6").Value = "2023-07-11 09:00:04"
$ws.Range("F6").Value = "2023-07-11 09:23:32"
$ws.Range("G6").Value = "00:23:28"
